$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.844999999999999
$ws.Range("B4").Value = 6.611999999999999
$ws.Range("C6").Value = -12.52
$ws.Range("B7").Value = 6.411
$ws.Range("C7").Value = -12.918
$ws.Range("B8").Value = 5.823
$ws.Range("C8").Value = -12.241
$ws.Range("A11").Value = -21.643
$ws.Range("D11").Value = -8.334
$ws.Range("A12").Value = -21.456
$ws.Range("B12").Value = 6.695
$ws.Range("B14").Value = 6.882000000000001
$ws.Range("D14").Value = -7.764
$ws.Range("A15").Value = -21.192
$ws.Range("C19").Value = -12.69
$ws.Range("D19").Value = -7.742999999999999
$ws.Range("C21").Value = -12.588
$ws.Range("D21").Value = -7.528
$ws.Range("B22").Value = 6.686
$ws.Range("C24").Value = -12.255
$ws.Range("C25").Value = -12.69
